$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template values shared across the new rows (same pattern as rows 4-79)
$bVal = 116.4121952
$cVal = 0.00170247
$dVal = 0.008850780000000001
$eVal = 0.06933635
$fVal = 12792.90181321
$gVal = 465.80531254
$hVal = 0.24
$iVal = 1.7904431
$jVal = 485.38834923

$dates = @(45636, 45637, 45638)
$startRow = 80
$templateRow = 79

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i

    # Copy the formatting (style) of the template row down to the new row first,
    # so the new cells inherit the same cell style (border/font/alignment/number format).
    $ws.Range("A$templateRow`:J$templateRow").Copy() | Out-Null
    $ws.Range("A$r`:J$r").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 3).Value = $cVal
    $ws.Cells.Item($r, 4).Value = $dVal
    $ws.Cells.Item($r, 5).Value = $eVal
    $ws.Cells.Item($r, 6).Value = $fVal
    $ws.Cells.Item($r, 7).Value = $gVal
    $ws.Cells.Item($r, 8).Value = $hVal
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

$excel.CutCopyMode = 0
